$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Update Price (D) and Volume(1h) (E) columns for rows 2-41
Set-TextValue "D2" "43.996.87"
Set-TextValue "E2" "  -0.95%  "
Set-TextValue "D3" "2.198.30"
Set-TextValue "E3" "  -2.17%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.10%  "
Set-TextValue "D5" "295.63"
Set-TextValue "E5" "  -4.16%  "
Set-TextValue "D6" "89.69"
Set-TextValue "E6" "  -6.15%  "
Set-TextValue "D7" "0.564"
Set-TextValue "E7" "  -1.46%  "
Set-TextValue "E8" "  -0.09%  "
Set-TextValue "D9" "0.486"
Set-TextValue "E9" "  -8.03%  "
Set-TextValue "D10" "32.41"
Set-TextValue "E10" "  -8.01%  "
Set-TextValue "D11" "0.0774"
Set-TextValue "E11" "  -4.79%  "
Set-TextValue "E12" "  -1.32%  "
Set-TextValue "D13" "6.87"
Set-TextValue "E13" "  -5.24%  "
Set-TextValue "D14" "2.528.25"
Set-TextValue "E14" "  -2.32%  "
Set-TextValue "D15" "2.261.29"
Set-TextValue "E15" "  -4.19%  "
Set-TextValue "D16" "13.25"
Set-TextValue "E16" "  -3.60%  "
Set-TextValue "D17" "0.778"
Set-TextValue "E17" "  -7.73%  "
Set-TextValue "D18" "43.654.64"
Set-TextValue "E18" "  -1.04%  "
Set-TextValue "D19" "0.0₃0893"
Set-TextValue "E19" "  -7.86%  "
Set-TextValue "D20" "5.84"
Set-TextValue "E20" "  -9.04%  "
Set-TextValue "D21" "10.96"
Set-TextValue "E21" "  -11.64%  "
Set-TextValue "D22" "63.34"
Set-TextValue "E22" "  -3.95%  "
Set-TextValue "D23" "232.71"
Set-TextValue "E23" "  -2.18%  "
Set-TextValue "D24" "2.77"
Set-TextValue "E24" "  -13.67%  "
Set-TextValue "E25" "  +0.71%  "
Set-TextValue "D26" "1.84"
Set-TextValue "E26" "  -8.73%  "
Set-TextValue "D27" "2.24"
Set-TextValue "E27" "  +0.53%  "
Set-TextValue "D28" "36.62"
Set-TextValue "E28" "  -5.57%  "
Set-TextValue "D29" "9.29"
Set-TextValue "E29" "  -6.14%  "
Set-TextValue "D30" "19.32"
Set-TextValue "E30" "  -4.06%  "
Set-TextValue "D31" "148.48"
Set-TextValue "E31" "  -3.63%  "
Set-TextValue "D32" "5.33"
Set-TextValue "E32" "  -10.95%  "
Set-TextValue "D33" "2.52"
Set-TextValue "E33" "  -5.02%  "
Set-TextValue "D34" "0.0739"
Set-TextValue "E34" "  -8.07%  "
Set-TextValue "D35" "0.116"
Set-TextValue "E35" "  -3.98%  "
Set-TextValue "D36" "2.85"
Set-TextValue "E36" "  -9.41%  "
Set-TextValue "E37" "  -6.12%  "
Set-TextValue "D38" "1.65"
Set-TextValue "E38" "  -8.21%  "
Set-TextValue "D39" "0.0287"
Set-TextValue "E39" "  -5.75%  "
Set-TextValue "D40" "3.53"
Set-TextValue "E40" "  -8.07%  "
Set-TextValue "D41" "3.09"
Set-TextValue "E41" "  -12.23%  "

# Rows 42-51: coin list reshuffled (row order changed, new coin inserted, one dropped)
Set-TextValue "B42" "Celestia"
Set-TextValue "C42" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D42" "13.15"
Set-TextValue "E42" "  -10.62%  "
Set-TextValue "B43" "FirstDigitalUSD"
Set-TextValue "C43" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D43" "1.01"
Set-TextValue "E43" "  -0.24%  "
Set-TextValue "B44" "Maker"
Set-TextValue "C44" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D44" "1.800.84"
Set-TextValue "E44" "  +2.99%  "
Set-TextValue "B45" "Stacks"
Set-TextValue "C45" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D45" "1.68"
Set-TextValue "E45" "  +4.51%  "
Set-TextValue "B46" "EnergySwap"
Set-TextValue "C46" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "14.18"
Set-TextValue "E46" "  +8.78%  "
Set-TextValue "B47" "HuobiToken"
Set-TextValue "C47" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D47" "2.82"
Set-TextValue "E47" "  +11.36%  "
Set-TextValue "B48" "Algorand"
Set-TextValue "C48" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D48" "0.175"
Set-TextValue "E48" "  -10.15%  "
Set-TextValue "B49" "BitcoinSV"
Set-TextValue "C49" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D49" "72.84"
Set-TextValue "E49" "  -9.67%  "
Set-TextValue "B50" "Aave"
Set-TextValue "C50" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D50" "92.27"
Set-TextValue "E50" "  -7.96%  "
Set-TextValue "B51" "ordi"
Set-TextValue "C51" "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue "D51" "65.35"
Set-TextValue "E51" "  -7.88%  "
